# Update teacher info in the economic_game workbook:
#  - add two new columns: "subject" and "company_money"
#  - reorder the "new student" block (Сагайдак Полина / Сагайдак Илья /
#    Симошин Михаил / Васильева Татьяна) up to rows 9-12
#  - move the two "teacher" rows (Валентинова / Игнатов) to the bottom
#    (rows 13-14), clear their now-unused grade/email cells, and give
#    them a subject + company_money value instead

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- header row: two new trailing columns -------------------------------
$ws.Range("F1").Value = "subject"
$ws.Range("G1").Value = "company_money"

# --- rows 9-12: student records (moved up, content unchanged) ----------
$ws.Range("A9").Value  = "Сагайдак"
$ws.Range("B9").Value  = "Полина"
$ws.Range("C9").Value  = "Тарасовна"
$ws.Range("D9").Value  = "12Б"
$ws.Range("E9").Value  = "polina.sagaidac@mail.ru"

$ws.Range("A10").Value = "Сагайдак"
$ws.Range("B10").Value = "Илья"
$ws.Range("C10").Value = "Тарасович"
$ws.Range("D10").Value = "12А"
$ws.Range("E10").Value = "m.s.v.inkognito@yandex.ru"

$ws.Range("A11").Value = "Симошин"
$ws.Range("B11").Value = "Михаил"
$ws.Range("C11").Value = "Николаевич"
$ws.Range("D11").Value = "9А"
$ws.Range("E11").Value = "mih_z8@mail.ru"

$ws.Range("A12").Value = "Васильева"
$ws.Range("B12").Value = "Татьяна"
$ws.Range("C12").Value = "Павловна"
$ws.Range("D12").Value = "11А"
$ws.Range("E12").Value = "sdelorec@yandex.ru"

# --- rows 13-14: teacher rows moved to the bottom -----------------------
# Content swaps with what used to be in rows 9 & 10; grade/email are no
# longer applicable for a teacher row, so those cells are removed
# entirely (not merely blanked), and subject/company_money are filled in.
$ws.Range("A13").Value = "Валентинова"
$ws.Range("B13").Value = "Валентина"
$ws.Range("C13").Value = "Валентиновна"
$ws.Range("D13").Clear()
$ws.Range("E13").Clear()
$ws.Range("F13").Value = "Кругосвет"
$ws.Range("G13").Value = 1000

$ws.Range("A14").Value = "Игнатов"
$ws.Range("B14").Value = "Игнат"
$ws.Range("C14").Value = "Игнатьевич"
$ws.Range("D14").Clear()
$ws.Range("E14").Clear()
$ws.Range("F14").Value = "Цифромир"
$ws.Range("G14").Value = 2000

# --- match the existing body style (s="1") on every newly-created cell -
# Copy + PasteSpecial(xlPasteFormats) reuses the existing cellXf instead
# of allocating a new (duplicate) one the way Range.Style assignment
# would; it has to run AFTER the Value assignments above or the engine
# drops the pasted format as soon as a fresh value is written in.
$ws.Range("E1").Copy()
$ws.Range("F1:G1").PasteSpecial(-4122)
$ws.Range("F13:G13").PasteSpecial(-4122)
$ws.Range("F14:G14").PasteSpecial(-4122)
$ws.Range("D9:E10").PasteSpecial(-4122)
